# Update "想去人数" (F column) values across sheets, per the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1)
$ws1.Range("F2").Value = 182
$ws1.Range("F3").Value = 399
$ws1.Range("F4").Value = 1117
$ws1.Range("F8").Value = 1053
$ws1.Range("F9").Value = 523
$ws1.Range("F10").Value = 314
$ws1.Range("F13").Value = 299
$ws1.Range("F14").Value = 344
$ws1.Range("F15").Value = 21
$ws1.Range("F17").Value = 387
$ws1.Range("F18").Value = 431
$ws1.Range("F19").Value = 5482
$ws1.Range("F21").Value = 1535
$ws1.Range("F22").Value = 353
$ws1.Range("F23").Value = 4628
$ws1.Range("F24").Value = 4628
$ws1.Range("F27").Value = 1467
$ws1.Range("F30").Value = 638
$ws1.Range("F31").Value = 21
$ws1.Range("F33").Value = 3784

# Sheet "本地生活" (sheet3)
$ws3.Range("F2").Value = 9369

# Sheet "全部类型" (sheet4) - combined view mirroring the other sheets
$ws4.Range("F2").Value = 9369
$ws4.Range("F5").Value = 182
$ws4.Range("F6").Value = 399
$ws4.Range("F7").Value = 1117
$ws4.Range("F11").Value = 1053
$ws4.Range("F12").Value = 523
$ws4.Range("F13").Value = 314
$ws4.Range("F16").Value = 299
$ws4.Range("F17").Value = 344
$ws4.Range("F18").Value = 21
$ws4.Range("F23").Value = 387
$ws4.Range("F24").Value = 431
$ws4.Range("F25").Value = 5482
$ws4.Range("F27").Value = 1535
$ws4.Range("F30").Value = 353
$ws4.Range("F32").Value = 4628
$ws4.Range("F33").Value = 4628
$ws4.Range("F36").Value = 1467
$ws4.Range("F39").Value = 638
$ws4.Range("F40").Value = 21
$ws4.Range("F47").Value = 3784

$wb.Save()
